$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-detected as numbers by Excel, so they stay literal text like the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "28.920.03"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "1.833.63"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "244.75"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "0.6901"
$ws.Range("E6").Value = "  -2.09%  "
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.07687"
$ws.Range("E8").Value = "  -2.95%  "
$ws.Range("D10").Value = "23.42"
$ws.Range("E10").Value = "  -4.48%  "
$ws.Range("D11").Value = "0.07785"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "1.824.78"
$ws.Range("E12").Value = "  -3.95%  "
$ws.Range("D13").Value = "5.082"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").Value = "90.41"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "0.6805"
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("D16").Value = "6.438"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "0.000008289"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "28.902.19"
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("D19").Value = "242.77"
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("D20").Value = "2.076.47"
$ws.Range("E20").Value = "  -3.62%  "
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "7.477"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "0.1474"
$ws.Range("E25").Value = "  -4.97%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "162.03"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "8.812"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "18.20"
$ws.Range("E28").Value = "  -3.28%  "
$ws.Range("D29").Value = "1.543"
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("D30").Value = "4.210"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("D33").Value = "0.05111"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("D34").Value = "0.7644"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "0.01849"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "1.218.74"
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("D40").Value = "2.698"
$ws.Range("D41").Value = "0.9392"
$ws.Range("D42").Value = "108.11"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "5.686"
$ws.Range("E44").Value = "  -5.80%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.00000000122"
$ws.Range("E45").Value = "  -3.22%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.5164"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.520"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").Value = "1.976.42"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("D49").Value = "64.10"
$ws.Range("E49").Value = "  -9.55%  "
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("E51").Value = "  -2.74%  "
